# Mark "Tested" (H) and "Done" (I) columns as "yes" for the backlog items
# that just got finished, matching the same visual style ("yes" look) that
# is already used on the "Logic"/"Design" (F/G) columns of each of those
# rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 6, 15, 16, 17, 18)

foreach ($r in $rows) {
    # Copy the formatting already used for "yes" on this row (column G)
    # onto the H and I cells, then set their value to "yes".
    $ws.Range("G${r}").Copy()
    $ws.Range("H${r}:I${r}").PasteSpecial(-4122)

    $ws.Range("H${r}").Value = "yes"
    $ws.Range("I${r}").Value = "yes"
}

# Update the active selection to match where the user ended up.
$ws.Range("H15:I18").Select()
